# Add 5 new literature-review rows about "TRM" (Tasa Representativa del
# Mercado) to the "Variables" sheet / Tabla1, matching the author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$lo = $ws.ListObjects.Item(1)

# Expand the table from A1:D8 to A1:D13 so the new rows become part of
# Tabla1 (autofilter + styling follow the resize).
$lo.Resize($ws.Range("A1:D13"))

# Copy the formatting of the last existing data row (row 8) down onto each
# of the five freshly-added rows (one row at a time, so the 4-column
# source tiles exactly onto each 4-column destination) before the cell
# values are overwritten.
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 9: Comparación y evaluación de pronosticos de la TRM
$ws.Range("A9").Value = "TRM"
$ws.Range("B9").Value = "Paper"
$ws.Range("C9").Value = "Comparación y evaluación de pronosticos de la Tasa Representativa del Mercado TRM"
$ws.Range("D9").Value = "https://bit.ly/3LdaEOM"
$ws.Rows.Item(9).RowHeight = 28.5

# Row 10: Modelación y pronóstico de la TRM...
$ws.Range("A10").Value = "TRM"
$ws.Range("B10").Value = "Paper"
$ws.Range("C10").Value = "Modelación y pronóstico de la TRM a partir de un modelo de saltos de difusión, un modelo de Black and Scholes y un modelo ARIMA, así como la comparación de los resultados de pronóstico entre ellos"
$ws.Range("D10").Value = "https://bit.ly/3qEZfzA"
$ws.Rows.Item(10).RowHeight = 57

# Row 11: EFECTOS DE LOS CHOQUES DE LOS PRECIOS EN EL WTI...
$ws.Range("A11").Value = "TRM y Petroleo WTI"
$ws.Range("B11").Value = "Paper"
$ws.Range("C11").Value = "EFECTOS DE LOS CHOQUES DE LOS PRECIOS EN EL WTI SOBRE LA TASA DE CAMBIO (TRM) EN COLOMBIA"
$ws.Range("D11").Value = "https://bit.ly/3uznfW3"
$ws.Rows.Item(11).RowHeight = 42.75

# Row 12: Efecto de las variables macroeconómicas globales y locales...
$ws.Range("A12").Value = "TRM "
$ws.Range("B12").Value = "Paper"
$ws.Range("C12").Value = "Efecto de las variables macroeconómicas globales y locales sobre el comportamiento de los futuros de la TRM en Colombia"
$ws.Range("D12").Value = "https://bit.ly/36Nkrwc"
$ws.Rows.Item(12).RowHeight = 42.75

# Row 13: PROYECCIÓN DE LA TASA DE CAMBIO DE COLOMBIA...
# (column C on this row was pasted in from elsewhere without the usual
# table border/font, just a plain wrap-text cell, so clear it first)
$ws.Range("A13").Value = "TRM"
$ws.Range("B13").Value = "Ppaper"
$ws.Range("C13").ClearFormats()
$ws.Range("C13").Value = "PROYECCIÓN DE LA TASA DE CAMBIO DE COLOMBIA BAJO CONDICIONES DE PPA: EVIDENCIA EMPÍRICA USANDO VAR"
$ws.Range("C13").WrapText = $true
$ws.Range("D13").Value = "https://bit.ly/36pqjfv"
$ws.Rows.Item(13).RowHeight = 30

# Only the first of the five new rows keeps a real hyperlink object (matches
# the uploaded workbook); the remaining new "Fuente" cells stay plain text.
$ws.Hyperlinks.Add($ws.Range("D9"), "https://bit.ly/3LdaEOM") | Out-Null
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
